$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename "SM" -> "SW" and "SM_Emergency" -> "SW_Emergency" throughout the
# "main" sheet (board identifier renamed).
$ws1.Range("D4").Value = "SW"
$ws1.Range("C6").Value = "SW"
$ws1.Range("D7").Value = "SW"
$ws1.Range("D11").Value = "SW"
$ws1.Range("E3").Value = "SW_Emergency"

# Remove the trailing empty rows on the "identifiers" sheet.
$ws2.Range("A13:D16").Delete() | Out-Null

# Update selections / active sheet to match the saved view state.
$ws1.Range("D11").Select() | Out-Null
$ws2.Activate()
$ws2.Range("A13:D16").Select() | Out-Null
